$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price values are plain decimals (e.g. "384.67") would be
# auto-converted to numbers by COMs smart-typing if assigned directly, but
# the source data keeps these as text (inline strings). We write the text via a
# literal formula, then convert that cell to a plain value in place (copy + paste
# special values-only) so the final cell is a plain text value, not a formula and
# not a number, and no new number-format style gets attached to the cell.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $escaped = $val.Replace("""", """""")
    $c.Formula = "=""" + $escaped + """"
    $c.Copy() | Out-Null
    $c.PasteSpecial(-4163) | Out-Null
}

$ws.Range("D2").Value = "51.677.21"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "3.058.52"
$ws.Range("E3").Value = "  +3.39%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "384.67"
$ws.Range("E5").Value = "  +1.33%  "
Set-TextValue "D6" "103.53"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  +0.03%  "
Set-TextValue "D9" "0.586"
$ws.Range("E9").Value = "  -1.03%  "
Set-TextValue "D10" "37.16"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").Value = "3.547.30"
$ws.Range("E13").Value = "  +3.65%  "
Set-TextValue "D14" "18.72"
$ws.Range("E14").Value = "  +2.43%  "
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "3.057.81"
$ws.Range("E16").Value = "  +3.30%  "
Set-TextValue "D17" "0.978"
$ws.Range("E17").Value = "  -1.72%  "
Set-TextValue "D18" "10.52"
$ws.Range("E18").Value = "  -5.44%  "
$ws.Range("D19").Value = "51.734.81"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("E20").Value = "  +0.21%  "
Set-TextValue "D21" "12.48"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "0.0₃0966"
$ws.Range("E22").Value = "  +0.51%  "
Set-TextValue "D23" "70.31"
$ws.Range("E23").Value = "  +0.07%  "
Set-TextValue "D24" "269.04"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("E25").Value = "  -1.23%  "
Set-TextValue "D26" "8.49"
$ws.Range("E26").Value = "  +8.34%  "
Set-TextValue "D27" "27.02"
$ws.Range("E27").Value = "  +4.41%  "
Set-TextValue "D28" "0.172"
$ws.Range("E28").Value = "  +5.53%  "
Set-TextValue "D29" "7.29"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("E32").Value = "  +0.37%  "
Set-TextValue "D33" "34.53"
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("E34").Value = "  +0.40%  "
Set-TextValue "D35" "50.47"
$ws.Range("E35").Value = "  -1.30%  "
Set-TextValue "D36" "0.0446"
$ws.Range("E36").Value = "  +2.65%  "
$ws.Range("E37").Value = "  -0.07%  "
Set-TextValue "D38" "3.39"
$ws.Range("E38").Value = "  +4.47%  "
Set-TextValue "D39" "0.291"
$ws.Range("E39").Value = "  +7.06%  "
Set-TextValue "D40" "17.10"
$ws.Range("E40").Value = "  +4.06%  "
$ws.Range("E41").Value = "  +2.95%  "
Set-TextValue "D42" "128.62"
$ws.Range("E42").Value = "  +2.75%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("E44").Value = "  +1.64%  "
Set-TextValue "D45" "3.74"
$ws.Range("E45").Value = "  +5.38%  "
Set-TextValue "D46" "22.03"
$ws.Range("E46").Value = "  +2.53%  "
$ws.Range("E47").Value = "  +6.22%  "
$ws.Range("E48").Value = "  +3.34%  "
$ws.Range("D49").Value = "2.048.81"
$ws.Range("E49").Value = "  +0.53%  "
$ws.Range("D50").Value = "3.365.87"
$ws.Range("E50").Value = "  +3.53%  "
$ws.Range("E51").Value = "  +7.33%  "

$excel.CutCopyMode = 0

